$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 54713
$ws.Range("F4").Value = 1336
$ws.Range("F5").Value = 377
$ws.Range("F7").Value = 879
$ws.Range("F8").Value = 749
$ws.Range("F9").Value = 394
$ws.Range("F10").Value = 3052
$ws.Range("F11").Value = 898
$ws.Range("F12").Value = 5221
$ws.Range("F14").Value = 1006
$ws.Range("F16").Value = 844
$ws.Range("F19").Value = 1271
$ws.Range("F22").Value = 175
$ws.Range("F23").Value = 361
$ws.Range("F24").Value = 24
$ws.Range("F25").Value = 37
$ws.Range("F29").Value = 4999
$ws.Range("F31").Value = 4935
$ws.Range("F32").Value = 8930
$ws.Range("F35").Value = 134
$ws.Range("F36").Value = 219
$ws.Range("F37").Value = 423
$ws.Range("F40").Value = 4203
$ws.Range("F41").Value = 240

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 55

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 571

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1336
$ws.Range("F5").Value = 879
$ws.Range("F6").Value = 749
$ws.Range("F7").Value = 394
$ws.Range("F8").Value = 3052
$ws.Range("F9").Value = 898
$ws.Range("F14").Value = 1006
$ws.Range("F15").Value = 55
$ws.Range("F16").Value = 844
$ws.Range("F19").Value = 1271
$ws.Range("F22").Value = 175
$ws.Range("F24").Value = 361
$ws.Range("F25").Value = 24
$ws.Range("F26").Value = 37
$ws.Range("F28").Value = 4999
$ws.Range("F30").Value = 8930
$ws.Range("F34").Value = 134
$ws.Range("F35").Value = 219
$ws.Range("F36").Value = 423
$ws.Range("F41").Value = 4203
$ws.Range("F48").Value = 240
